# Adding code for create new document functionality
# Adds a new "Docs" worksheet at the end of the workbook, populates it with
# two sample documents, formats the header row, and makes it the active sheet
# (mirroring what the "Calls" sheet loses when it stops being the active tab).

$wb = $excel.ActiveWorkbook

# --- Remember the previously-active sheet (Calls) so we can update its
#     selection to a full-row selection once it's no longer the active tab.
$calls = $wb.Worksheets.Item("Calls")

# --- Create the new "Docs" sheet after the last existing sheet ("Calls") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$docs = $wb.Worksheets.Add($null, $lastSheet)
$docs.Name = "Docs"

# --- Header row ---
$headers = @("title", "description", "version", "contact", "client", "prospect", "task", "case", "tags")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $docs.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$docs.Range("A1:I1").Interior.Color = 65535

# --- Data rows (version values are quote-prefixed text, like "0.1"/"1.9") ---
$docs.Range("A2").Value = "Doc1"
$docs.Range("B2").Value = "Document 1"
$docs.Range("C2").Value = "'0.1"
$docs.Range("D2").Value = "abcd"
$docs.Range("E2").Value = "efgh"
$docs.Range("F2").Value = "ijkl"
$docs.Range("G2").Value = "mnop"
$docs.Range("H2").Value = "qrst"
$docs.Range("I2").Value = "uvwxyz"

$docs.Range("A3").Value = "Doc2"
$docs.Range("B3").Value = "Document 2"
$docs.Range("C3").Value = "'1.9"
$docs.Range("D3").Value = "dcba"
$docs.Range("E3").Value = "hgfe"
$docs.Range("F3").Value = "lkji"
$docs.Range("G3").Value = "ponm"
$docs.Range("H3").Value = "tsrq"
$docs.Range("I3").Value = "zyxwvy"

# --- Column B ("description") is wide enough to fit its content ---
$docs.Range("B1:B3").EntireColumn.AutoFit()

# --- Update the old active sheet's (Calls) selection to a full header-row
#     selection, then activate the new Docs sheet so it becomes the tab
#     shown when the workbook is opened. ---
$calls.Rows(1).EntireRow.Select()
$docs.Range("J3").Select()
$docs.Activate()
